$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 2433.111
$ws.Range("J9").Value = 2128.5715
$ws.Range("L9").Value = 2128.5715
$ws.Range("N9").Value = -2466.5715

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 574.25
$ws.Range("I19").Value = 650
$ws.Range("J19").Value = 498.5
$ws.Range("K19").Value = 650
$ws.Range("L19").Value = 498.5
$ws.Range("M19").Value = -475
$ws.Range("N19").Value = -848.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 2164.1667
$ws.Range("J39").Value = 3330
$ws.Range("L39").Value = 9990
$ws.Range("N39").Value = -10582

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H42").Value = 901.2
$ws.Range("I42").Value = 998.5
$ws.Range("J42").Value = 836.3333
$ws.Range("K42").Value = 2995.5
$ws.Range("L42").Value = 2508.9999
$ws.Range("M42").Value = -2765.5
$ws.Range("N42").Value = -2968.9999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H46").Value = 2999.5
$ws.Range("I46").Value = 2999.5
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 8998.5
$ws.Range("L46").Value = 0
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -8879.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H59").Value = 4166.6665
$ws.Range("I59").Value = 2500
$ws.Range("K59").Value = 7500
$ws.Range("M59").Value = -6943

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H60").Value = 2999.5
$ws.Range("I60").Value = 2999.5
$ws.Range("J60").Value = 0
$ws.Range("K60").Value = 8998.5
$ws.Range("L60").Value = 0
$ws.Range("M60").ClearContents()
$ws.Range("N60").Value = -8514.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 3100
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 3100
$ws.Range("K70").Value = 0
$ws.Range("L70").ClearContents()
$ws.Range("M70").Value = 9300
$ws.Range("N70").Value = -9840

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H73").Value = 3100
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 3100
$ws.Range("K73").Value = 0
$ws.Range("L73").ClearContents()
$ws.Range("M73").Value = 9300
$ws.Range("N73").Value = -11172

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 2980

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H83").Value = 2980

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 996.25
$ws.Range("J98").Value = 990
$ws.Range("L98").Value = 990
$ws.Range("N98").Value = -3986

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 996.25
$ws.Range("J122").Value = 990
$ws.Range("L122").Value = 2970
$ws.Range("N122").Value = -7870

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 5790.1875
$ws.Range("I132").Value = 5297
$ws.Range("J132").Value = 6086.1
$ws.Range("K132").Value = 15891
$ws.Range("L132").Value = 18258.3
$ws.Range("M132").Value = -13361
$ws.Range("N132").Value = -23318.3

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H50").Value = 43750
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 43750
$ws.Range("K50").Value = 0
$ws.Range("L50").ClearContents()
$ws.Range("M50").Value = 43750
$ws.Range("N50").Value = -45178

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H135").Value = 60001
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 60001
$ws.Range("K135").Value = 0
$ws.Range("L135").ClearContents()
$ws.Range("M135").Value = 60001
$ws.Range("N135").Value = -70141

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1417.5
$ws.Range("I20").Value = 1764.4445
$ws.Range("J20").Value = 971.4286
$ws.Range("K20").Value = 1764.4445
$ws.Range("L20").Value = 971.4286
$ws.Range("M20").Value = -1517.4445
$ws.Range("N20").Value = -1465.4286

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 431.66666
$ws.Range("I22").Value = 523.3333
$ws.Range("J22").Value = 340
$ws.Range("K22").Value = 523.3333
$ws.Range("L22").Value = 340
$ws.Range("M22").Value = -350.3333
$ws.Range("N22").Value = -686

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H36").Value = 2500
$ws.Range("I36").Value = 2500
$ws.Range("K36").Value = 2500
$ws.Range("M36").Value = -1966

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H62").Value = 50000
$ws.Range("J62").Value = 50000
$ws.Range("L62").Value = 50000
$ws.Range("N62").Value = -51372

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H65").Value = 50000
$ws.Range("J65").Value = 50000
$ws.Range("L65").Value = 150000
$ws.Range("N65").Value = -156864

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H92").Value = 25000
$ws.Range("J92").Value = 25000
$ws.Range("L92").Value = 25000
$ws.Range("N92").Value = -29992

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 1608.2142
$ws.Range("I7").Value = 1210.9
$ws.Range("J7").Value = 2601.5
$ws.Range("K7").Value = 1210.9
$ws.Range("L7").Value = 2601.5
$ws.Range("M7").Value = -1097.9
$ws.Range("N7").Value = -2827.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H35").Value = 756.25
$ws.Range("I35").Value = 756.25
$ws.Range("K35").Value = 756.25
$ws.Range("M35").Value = -462.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H39").Value = 2536.5
$ws.Range("I39").Value = 2536.5
$ws.Range("K39").Value = 2536.5
$ws.Range("M39").Value = -2145.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H49").Value = 2536.5
$ws.Range("I49").Value = 2536.5
$ws.Range("K49").Value = 2536.5
$ws.Range("M49").Value = -2354.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H55").Value = 7272.727
$ws.Range("I55").Value = 7000
$ws.Range("K55").Value = 7000
$ws.Range("M55").Value = -6685

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H96").Value = 23504.666
$ws.Range("J96").Value = 23504.666
$ws.Range("L96").Value = 23504.666
$ws.Range("N96").Value = -28996.666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H47").Value = 74
$ws.Range("I47").Value = 74
$ws.Range("K47").Value = 222
$ws.Range("M47").Value = 209

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H50").Value = 366.42856
$ws.Range("I50").Value = 176.5
$ws.Range("K50").Value = 529.5
$ws.Range("M50").Value = -48.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H53").Value = 366.42856
$ws.Range("I53").Value = 176.5
$ws.Range("K53").Value = 529.5
$ws.Range("M53").Value = -48.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H60").Value = 4047.5
$ws.Range("J60").Value = 8000
$ws.Range("L60").Value = 24000
$ws.Range("N60").Value = -24502

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("M29").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2172.3125
$ws.Range("I122").Value = 1482.6428
$ws.Range("J122").Value = 7000
$ws.Range("K122").Value = 4447.928400000001
$ws.Range("L122").Value = 21000
$ws.Range("M122").Value = -1997.928400000001
$ws.Range("N122").Value = -25900

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 5213.8
$ws.Range("I126").Value = 5505.5713
$ws.Range("K126").Value = 16516.7139
$ws.Range("M126").Value = -14046.7139

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 354.75
$ws.Range("I16").Value = 354.75
$ws.Range("K16").Value = 354.75
$ws.Range("M16").Value = -184.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1016
$ws.Range("J22").Value = 881.5
$ws.Range("L22").Value = 881.5
$ws.Range("N22").Value = -1471.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 1016
$ws.Range("J27").Value = 881.5
$ws.Range("L27").Value = 881.5
$ws.Range("N27").Value = -1095.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H31").Value = 4346
$ws.Range("J31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("N31").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H35").Value = 4799.4
$ws.Range("I35").Value = 374.25
$ws.Range("J35").Value = 22500
$ws.Range("K35").Value = 374.25
$ws.Range("L35").Value = 22500
$ws.Range("M35").Value = -38.25
$ws.Range("N35").Value = -23172

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H53").Value = 19761.5
$ws.Range("J53").Value = 13500
$ws.Range("L53").Value = 13500
$ws.Range("N53").Value = -14536

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H55").Value = 14675.75
$ws.Range("J55").Value = 19517.666
$ws.Range("L55").Value = 19517.666
$ws.Range("N55").Value = -20071.666
